{"js": "// Apply the four text corrections described by the diff:\n// 1. \"185000\" -> \"185 000\"\n// 2. \"cent quatre-vingt-cinq mille \" -> \"cent quatre-vingt-cinq \" (drop \"mille\")\n// 3. drop \" et quatre-vingt-seize\" from the FCFA-in-words parenthetical\n// 4. \"22 novembre 2024\" -> \"02 d\u00e9cembre 2024\"\n\nasync function replaceOnce(searchText, replacement, options) {\n  const results = context.document.body.search(searchText, options || { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nawait replaceOnce(\"185000\", \"185 000\", { matchCase: true });\n\nawait replaceOnce(\n  \"cent quatre-vingt-cinq mille \",\n  \"cent quatre-vingt-cinq \",\n  { matchCase: true }\n);\n\nawait replaceOnce(\n  \"8 624 827 FCFA (huit millions six cent vingt-quatre mille huit cent vingt-six et quatre-vingt-seize CFA)\",\n  \"8 624 827 FCFA (huit millions six cent vingt-quatre mille huit cent vingt-six CFA)\",\n  { matchCase: true }\n);\n\nawait replaceOnce(\"22 novembre 2024\", \"02 d\u00e9cembre 2024\", { matchCase: true });\n", "ps1": "# Apply the four text corrections described by the diff:\n# 1. \"185000\" -> \"185 000\"\n# 2. \"cent quatre-vingt-cinq mille \" -> \"cent quatre-vingt-cinq \" (drop \"mille\")\n# 3. drop \" et quatre-vingt-seize\" from the FCFA-in-words parenthetical\n# 4. \"22 novembre 2024\" -> \"02 d\u00e9cembre 2024\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $findText,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        $wdFindContinue,  # Wrap\n        $false,      # Format\n        $replaceText,     # ReplaceWith\n        $wdReplaceAll     # Replace\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-Once \"185000\" \"185 000\"\nReplace-Once \"cent quatre-vingt-cinq mille \" \"cent quatre-vingt-cinq \"\nReplace-Once \"8 624 827 FCFA (huit millions six cent vingt-quatre mille huit cent vingt-six et quatre-vingt-seize CFA)\" \"8 624 827 FCFA (huit millions six cent vingt-quatre mille huit cent vingt-six CFA)\"\nReplace-Once \"22 novembre 2024\" \"02 d\u00e9cembre 2024\"\n\n$d.Save()\n"}
